$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "uk" flag code to "gb" (this also triggers the shared-string
# cleanup/reorder that drops the now-unused "uk" entry and appends "gb").
$ws.Range("G4").Value = "gb"

# Update the percent values (column B) per the sort fix.
$ws.Range("B2").Value = 60
$ws.Range("B4").Value = 45
$ws.Range("B5").Value = 17
$ws.Range("B6").Value = 24
$ws.Range("B7").Value = 8
$ws.Range("B8").Value = 12
$ws.Range("B9").Value = 13
$ws.Range("B10").Value = 15
$ws.Range("B11").Value = 12
$ws.Range("B12").Value = 5
$ws.Range("B13").Value = 23
$ws.Range("B14").Value = 12
$ws.Range("B16").Value = 20
$ws.Range("B17").Value = 6
$ws.Range("B18").Value = 5
$ws.Range("B20").Value = 35
$ws.Range("B21").Value = 14
$ws.Range("B22").Value = 33
$ws.Range("B23").Value = 2
$ws.Range("B24").Value = 12

# Update the selection to match the saved view state.
[void]$ws.Range("B24").Select()
